# v3.0 update FCI 27/1/2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header date for column C, matching the existing style used by B1
# (bold, bordered, centered) rather than the plain default style.
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "13-01-2023"

# Final row order/content (rows 2..16), each with the existing (B) and new (C) values.
$data = @(
    @("Alpha Acciones", 66310.58, 66612.98),
    @("Alpha Mega", 380068.65, 378929.83),
    @("Alpha Recursos Naturales", 86338.23, 86191.5),
    @("Alpha planeam equil", 5069.15, 5104.76),
    @("Delta Acciones", 2273.61, 2148.14),
    @("Delta Recursos Naturales", 227152.92, 226257.5),
    @("Delta Select", 2476.34, 2877.33),
    @("Fima Acciones", 94080.63, 95068.33),
    @("Fima PB Acciones", 43792.85, 43890.4),
    @("HF Acciones Argentinas", 0, 0),
    @("HF Acciones Lideres", 77173, 76472.61),
    @("Supefondo RV", 152776.93, 156377.46),
    @("Toronto Trust Multimercado", 3930.93, 4090.78),
    @("avg", 87803.37, 88001.66),
    @("total", 1141443.82, 1144021.62)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $entry = $data[$i]

    # A-column keeps the bold/bordered/centered label style already present
    # on the sheet (copied from the existing A2 label cell).
    $ws.Range("A2").Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)
    $ws.Cells.Item($row, 1).Value = $entry[0]

    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
}
